$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.242.80"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "2.994.62"
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'" + "501.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.07%  "
$ws.Range("D6").Value = "'" + "137.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.25%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("E8").Value = "  -3.85%  "
$ws.Range("D9").Value = "'" + "7.29"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.31%  "
$ws.Range("D10").Value = "'" + "0.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.07%  "
$ws.Range("E11").Value = "  -4.15%  "
$ws.Range("D12").Value = "3.499.91"
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").Value = "'" + "26.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.23%  "
$ws.Range("D15").Value = "'" + "0.0000159"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.02%  "
$ws.Range("D16").Value = "57.261.50"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "'" + "6.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("D18").Value = "2.990.70"
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("D19").Value = "'" + "12.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.11%  "
$ws.Range("D20").Value = "'" + "7.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.81%  "
$ws.Range("D21").Value = "'" + "319.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.12%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "'" + "0.492"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("D25").Value = "'" + "63.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.00%  "
$ws.Range("D26").Value = "'" + "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -5.33%  "
$ws.Range("D28").Value = "0.0₃0890"
$ws.Range("E28").Value = "  -9.70%  "
$ws.Range("D29").Value = "'" + "6.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.36%  "
$ws.Range("D30").Value = "'" + "7.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.33%  "
$ws.Range("E31").Value = "  -4.62%  "
$ws.Range("E32").Value = "  -7.41%  "
$ws.Range("D33").Value = "'" + "20.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.08%  "
$ws.Range("D34").Value = "'" + "155.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").Value = "'" + "4.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("D36").Value = "'" + "5.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("D37").Value = "'" + "1.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.23%  "
$ws.Range("D38").Value = "'" + "24.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.77%  "
$ws.Range("D39").Value = "'" + "0.0662"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.12%  "
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").Value = "3.022.52"
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("D42").Value = "'" + "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "'" + "3.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("D44").Value = "'" + "0.646"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("D45").Value = "2.190.30"
$ws.Range("E45").Value = "  -6.34%  "
$ws.Range("E46").Value = "  -7.24%  "
$ws.Range("D47").Value = "'" + "5.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("D48").Value = "'" + "0.939"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.11%  "
$ws.Range("E49").Value = "  -5.02%  "
$ws.Range("D50").Value = "'" + "19.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.21%  "
$ws.Range("D51").Value = "'" + "1.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -12.82%  "
